{"js": "// Remove the trailing \"Requisitos\" section (its Heading2 title paragraph\n// and the ListBullet paragraph listing the weak-requirement courses)\n// from the end of the document body.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Requisitos\" heading paragraph by its style + text, and\n// delete it together with the paragraph that immediately follows it\n// (the bullet list of LOQ course requisites). Matching by content\n// rather than a hard-coded index keeps this robust to minor\n// differences in paragraph count.\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  if (para.style === \"Heading 2\" && para.text.trim() === \"Requisitos\") {\n    const next = i + 1 < items.length ? items[i + 1] : null;\n    if (next) {\n      next.delete();\n    }\n    para.delete();\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Requisitos\" section (its Heading 2 title paragraph\n# and the List Bullet paragraph listing the weak-requirement courses)\n# from the end of the document body.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $txt -eq \"Requisitos\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ge 1) {\n    # Delete the paragraph that follows the heading first (the bullet\n    # list of LOQ course requisites), then the heading itself, so the\n    # indices stay valid while deleting.\n    if ($targetIndex + 1 -le $d.Paragraphs.Count) {\n        $d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n    }\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n}\n"}
